# Horarios actualizados Linea 141 - 486
# Updates the scraped bus-schedule data across the 3 worksheets to the
# new scrape timestamp (05:42:22) and refreshed arrival rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "LP1912"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 05:42:22"
$ws1.Range("A3").Value = "Total filas: 18"

$ws1Data = @(
    @("05:42:22", "05:43", "14_ABASTO", 1, "LP1912"),
    @("05:42:22", "05:52", "17_ROMERO", 10, "LP1912"),
    @("05:42:22", "06:01", "16_SANTA ANA", 19, "LP1912"),
    @("05:42:22", "06:04", "10_OLMOS", 22, "LP1912"),
    @("05:42:22", "06:11", "215A_EL PATO", 29, "LP1912"),
    @("05:42:22", "06:24", "11_ETCHEVERRY", 42, "LP1912"),
    @("05:42:22", "06:27", "23_HERNANDEZ", 45, "LP1912"),
    @("05:42:22", "06:31", "16_SANTA ANA", 49, "LP1912"),
    @("05:42:22", "06:31", "17X38_ROMERO", 49, "LP1912"),
    @("05:42:22", "06:39", "225_C ROCA-H SUR", 57, "LP1912"),
    @("05:42:22", "06:54", "14_ABASTO", 72, "LP1912"),
    @("05:42:22", "07:04", "225_GOMEZ", 82, "LP1912"),
    @("05:42:22", "07:07", "215C_EL PATO", 85, "LP1912"),
    @("05:42:22", "07:14", "14X44_ABASTO", 92, "LP1912"),
    @("05:42:22", "07:21", "215A_EL PATO", 99, "LP1912"),
    @("05:42:22", "07:33", "23_HERNANDEZ", 111, "LP1912"),
    @("05:42:22", "07:36", "17X38_ROMERO", 114, "LP1912"),
    @("05:42:22", "07:37", "27_EL RETIRO", 115, "LP1912")
)

$r = 6
foreach ($row in $ws1Data) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $r++
}

# ---------------------------------------------------------------------------
# Sheet "LP1912-215"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 05:42:22"
$ws2.Range("A3").Value = "Total filas: 3"

$ws2Data = @(
    @("05:42:22", "06:11", "215A_EL PATO", 29, "LP1912"),
    @("05:42:22", "07:07", "215C_EL PATO", 85, "LP1912"),
    @("05:42:22", "07:21", "215A_EL PATO", 99, "LP1912")
)

$r = 6
foreach ($row in $ws2Data) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $r++
}

# ---------------------------------------------------------------------------
# Sheet "6203-6173"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 05:42:22"
$ws3.Range("A3").Value = "Total filas: 1"

# Bring in the header row (with its bold + bordered + centered style) from
# the first sheet instead of re-building the format from scratch.
$ws1.Range("A5:E5").Copy($ws3.Range("A5:E5"))

$ws3.Cells.Item(6, 1).Value = "05:42:22"
$ws3.Cells.Item(6, 2).Value = "07:27"
$ws3.Cells.Item(6, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(6, 4).Value = 105
$ws3.Cells.Item(6, 5).Value = "L6173"
